$d = $word.ActiveDocument

# Plain text replace - safe to use when the matched span does not cross any
# w:proofErr markers or runs carrying character formatting (keeps those
# intact because the match stays inside "plain" runs only).
function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Insert text immediately after a uniquely-identifying anchor string,
# without touching/merging anything that follows (used where a naive
# Find/Replace would otherwise swallow w:proofErr or formatted runs).
function Insert-After($anchor, $insertText) {
    $rng = $d.Content
    $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $pt = $d.Range($rng.End, $rng.End)
    $pt.InsertBefore($insertText)
}

# --- Paragraph starting "LcmsNet provides access to DMS functionality..." ---

Replace-Text "Lcmsnet to retrieve data from a desired data source and use it" "users to retrieve data from a desired data source(s) and use it"

Replace-Text "database/data storage system using the" "database/data storage system(s) using the"

# [Export(typeof(IDmsTools)] -> [Export(typeof(IDmsTools))]  (insert extra ")")
Insert-After "IDmsTools)" ")"

Replace-Text "MEF attribute and if desired [" "MEF attribute, and if desired  the ["

Replace-Text ")]. As of this writing, those are the only two pieces of metadata to use. The metadata" ")] MEF attribute. As of this writing, those are the only two pieces of metadata tracked. The metadata"

# --- Paragraph starting "After this, compile your library..." ---

Replace-Text "After this, compile your library and place the resulting library files in the" "After compiling your DMS library place the resulting library files in the"

Replace-Text "directory during or after install of LcmsNet.  Once complete, start LcmsNet, and select your tool in the Configuration window under the DMS tab, if it is not selected by " "directory during or after install of LcmsNet.  Once placed, start or restart LcmsNet and select your tool in the Configuration window under the DMS tab if it is not selected by "

Replace-Text "default." "default, LcmsNet will then be able to retrieve information from the desired data source."
